$wb = $excel.ActiveWorkbook

# ---- Update "PCA_VAR" sheet (sheet3) ----
$wsVar = $wb.Worksheets.Item("PCA_VAR")

$pcaVarValues = @(
    0.496, 0.141, 0.128, 0.077, 0.07000000000000001, 0.036, 0.018
)
for ($i = 0; $i -lt $pcaVarValues.Length; $i++) {
    $row = $i + 2
    $wsVar.Cells.Item($row, 2).Value = $pcaVarValues[$i]
}

# Remove old rows 9-13 (data no longer present after the update)
$wsVar.Range("A9:A13").EntireRow.Delete()

# ---- Update "PCA_Components" sheet (sheet4) ----
$wsComp = $wb.Worksheets.Item("PCA_Components")

$pcaCompValues = @(
    @(-0.015, 0.102, 0.033, 0.02, 0.035, -0.017, 0.6860000000000001, 0.243, 0.637, -0.074, -0.02, 0.075, 0.03, 0.014, 0, 0.194, 0.037),
    @(0.099, 0.266, 0.048, 0.058, -0.002, 0.013, 0.081, 0.823, -0.423, 0.056, 0.019, 0.064, 0.131, 0.043, 0.003, -0.113, 0.108),
    @(-0.161, -0.374, 0.01, 0.015, -0.005, -0.027, -0.047, 0.169, -0.226, -0.124, -0.04, -0.008999999999999999, -0.005, 0.008999999999999999, -0.022, 0.85, -0.109),
    @(-0.054, -0.421, -0.001, 0.005, -0.011, -0.005, 0.147, 0.192, -0.014, 0.005, -0.01, 0.002, -0.067, -0.027, -0.008, -0.332, -0.803),
    @(-0.097, -0.707, 0.008999999999999999, 0.01, 0.002, -0.049, 0.057, 0.146, 0.024, -0.232, -0.068, 0.023, 0.006, -0.005, -0.005, -0.317, 0.552),
    @(0.447, -0.283, 0.075, 0.054, 0.065, 0.171, -0.127, 0.002, 0.181, 0.484, 0.192, 0.096, 0.551, 0.19, 0.04, 0.08799999999999999, 0.003),
    @(0.191, 0.12, 0.102, 0.079, 0.079, -0.091, -0.168, -0.047, 0.061, -0.706, -0.266, 0.08699999999999999, 0.528, 0.051, 0.028, -0.031, -0.153)
)
for ($i = 0; $i -lt $pcaCompValues.Length; $i++) {
    $row = $i + 2
    $rowVals = $pcaCompValues[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 2  # Column B = 2
        $wsComp.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}

# Remove old rows 9-13 (data no longer present after the update)
$wsComp.Range("A9:A13").EntireRow.Delete()

Write-Output "Updated PCA_VAR and PCA_Components sheets successfully."
